# Add 4 new variable rows (V1/premotor myelin measures) into the data
# dictionary, inserted just above the existing "task_A_motion" row
# (old row 197), pushing the rest of the table down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new rows: insert 4 blank rows starting at row 197.
$ws.Rows("197:200").Insert()

# Row 197: premotor_A_myelin
$ws.Range("A197").Value = "premotor_A_myelin"
$ws.Range("B197").Value = "precentral gyrus from Harvard-Oxford probabilistic cortical strucutral atlas T1/T2 intensity ratio during scan A, calibrated using csf/skull"

# Row 198: premotor_B_myelin
$ws.Range("A198").Value = "premotor_B_myelin"
$ws.Range("B198").Value = "precentral gyrus from Harvard-Oxford probabilistic cortical strucutral atlas T1/T2 intensity ratio during scan B, calibrated using csf/skull"

# Row 199: V1_A_myelin
$ws.Range("A199").Value = "V1_A_myelin"
$ws.Range("B199").Value = "pericalcarine gyrus from Harvard-Oxford probabilistic cortical strucutral atlas T1/T2 intensity ratio during scan A, calibrated using csf/skull"

# Row 200: V1_B_myelin
$ws.Range("A200").Value = "V1_B_myelin"
$ws.Range("B200").Value = "pericalcarine gyrus from Harvard-Oxford probabilistic cortical strucutral atlas T1/T2 intensity ratio during scan B, calibrated using csf/skull"

# Match the saved view state: selection on B197, scrolled so row 181 is
# the top visible row.
$ws.Range("B197").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 181
$win.ScrollColumn = 1
